$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update "Riders" (column C) and "Average" (column D) values with new Madigan bike hours data
$ws.Range("C2").Value = 193
$ws.Range("D2").Value = 225.14

$ws.Range("C3").Value = 205
$ws.Range("D3").Value = 211.35

$ws.Range("C4").Value = 233
$ws.Range("D4").Value = 213.22

$ws.Range("C5").Value = 239
$ws.Range("D5").Value = 239.32

$ws.Range("C6").Value = 270
$ws.Range("D6").Value = 242.83

$ws.Range("C7").Value = 95
$ws.Range("D7").Value = 113.04

$ws.Range("C8").Value = 69
$ws.Range("D8").Value = 93.09

$wb.Save()
